# chore: update Sheets via scheduled runner
#
# Applies a scheduled-runner refresh of market/profit data across three
# leve-profit tables (ALC, CRP, LTW): some rows' price/profit columns
# (H:N) are cleared because no current market data is available, some
# rows get refreshed price/profit values, and a couple of rows have
# their computed profit values updated.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# ALC sheet: rows 125-141 lost all of their H:N (price/profit) data.
# ---------------------------------------------------------------------
$wsAlc = $wb.Worksheets.Item("ALC")
$wsAlc.Range("H125:N141").ClearContents()

# ---------------------------------------------------------------------
# CRP sheet: refreshed price data for rows 31 & 34, and newly-populated
# price/profit data for rows 129-141 (row 136 already had data and is
# untouched).
# ---------------------------------------------------------------------
$wsCrp = $wb.Worksheets.Item("CRP")

$wsCrp.Range("H31").Value = 1503
$wsCrp.Range("I31").Value = 1503
$wsCrp.Range("J31").Value = 0
$wsCrp.Range("K31").Value = 1503
$wsCrp.Range("L31").Value = 0
$wsCrp.Range("M31").Value = -1208
$wsCrp.Range("N31").ClearContents()

$wsCrp.Range("H34").Value = 1503
$wsCrp.Range("I34").Value = 1503
$wsCrp.Range("J34").Value = 0
$wsCrp.Range("K34").Value = 1503
$wsCrp.Range("L34").Value = 0
$wsCrp.Range("M34").Value = -1301
$wsCrp.Range("N34").ClearContents()

$wsCrp.Range("H129:L129").Value = 0

$wsCrp.Range("H130:L130").Value = 0

$wsCrp.Range("H131").Value = 30000
$wsCrp.Range("I131").Value = 0
$wsCrp.Range("J131").Value = 30000
$wsCrp.Range("K131").Value = 0
$wsCrp.Range("L131").Value = 30000
$wsCrp.Range("N131").Value = -40080

$wsCrp.Range("H132").Value = 5402
$wsCrp.Range("I132").Value = 5012
$wsCrp.Range("J132").Value = 5499.5
$wsCrp.Range("K132").Value = 15036
$wsCrp.Range("L132").Value = 16498.5
$wsCrp.Range("M132").Value = -12506
$wsCrp.Range("N132").Value = -21558.5

$wsCrp.Range("H133:L133").Value = 0

$wsCrp.Range("H134").Value = 3537.3333
$wsCrp.Range("I134").Value = 2112
$wsCrp.Range("J134").Value = 4250
$wsCrp.Range("K134").Value = 6336
$wsCrp.Range("L134").Value = 12750
$wsCrp.Range("M134").Value = -3801
$wsCrp.Range("N134").Value = -17820

$wsCrp.Range("H135:L135").Value = 0

$wsCrp.Range("H137:L137").Value = 0

$wsCrp.Range("H138:L138").Value = 0

$wsCrp.Range("H139:L139").Value = 0

$wsCrp.Range("H140:L140").Value = 0

$wsCrp.Range("H141:L141").Value = 0

# ---------------------------------------------------------------------
# LTW sheet: rows 124-141 (except 126, which already had no data and is
# unaffected) lost all of their H:N (price/profit) data.
# ---------------------------------------------------------------------
$wsLtw = $wb.Worksheets.Item("LTW")
$wsLtw.Range("H124:N125").ClearContents()
$wsLtw.Range("H127:N141").ClearContents()
